# plotlyjs_commands_example.pptx - "Consistently use double-quotes for all JSON snippets"
#
# Two kinds of changes:
#  1. The "datetimeFigureOut" date placeholders (slide master, all 11
#     layouts, and the notes master) get their cached text bumped from
#     5/19/18 to 5/22/18.
#  2. The Python-dict-looking JSON snippets on slide 1 get their single
#     quotes turned into double quotes so they read as valid JSON.

$p = $ppt.ActivePresentation
$q = [string][char]34

# ---------------------------------------------------------------------
# 1. Date placeholders: slide master, every custom layout, notes master
# ---------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "5/19/18") {
                $tr.Text = "5/22/18"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes
}

Update-DatePlaceholders $p.NotesMaster.Shapes

# ---------------------------------------------------------------------
# 2. Slide 1: swap the Python-style single quotes for double quotes in
#    the little code snippets drawn in the "Rectangle NN" textboxes.
#    Each edit rewrites exactly one existing run (same start/length as
#    the original run) so formatting / run boundaries stay untouched.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

# --- Shape 1 "Rectangle 51" ---
$tr = $slide.Shapes.Item(1).TextFrame.TextRange
$tr.Characters(1, 10).Text = '{' + $q + 'data' + $q + ': ['
$tr.Characters(12, 19).Text = '    {' + $q + 'type' + $q + ': ' + $q + 'bar' + $q + ','
$tr.Characters(32, 22).Text = '     ' + $q + 'y' + $q + ': [2, 3, 1]}],'
$tr.Characters(55, 14).Text = ' ' + $q + 'layout' + $q + ': {}}'

# --- Shape 2 "Rectangle 52" ---
$tr = $slide.Shapes.Item(2).TextFrame.TextRange
$tr.Characters(1, 2).Text = '{' + $q
$tr.Characters(14, 11).Text = $q + ': [-1, 3]}'

# --- Shape 3 "Rectangle 53" ---
$tr = $slide.Shapes.Item(3).TextFrame.TextRange
$tr.Characters(1, 19).Text = '{' + $q + 'type' + $q + ': ' + $q + 'scatter' + $q + ','
$tr.Characters(21, 20).Text = '     ' + $q + 'y' + $q + ': [3, 1, 2]}'

# --- Shape 4 "Rectangle 54" ---
$tr = $slide.Shapes.Item(4).TextFrame.TextRange
$tr.Characters(1, 30).Text = '{' + $q + 'data' + $q + ': {' + $q + 'name' + $q + ': [' + $q + 'A' + $q + ', ' + $q + 'B' + $q + '], '
# The run "          'marker.size': " (chars 32-56) becomes three runs:
#   "          \""  +  "marker.size"  +  "\": "
$tr.Characters(32, 11).Text = '          ' + $q
$tr.Characters(43, 11).Text = 'marker.size'
$tr.Characters(54, 3).Text = $q + ': '
$tr.Characters(90, 13).Text = ' ' + $q + 'layout' + $q + ': {' + $q
$tr.Characters(117, 2).Text = $q + ':'

# --- Shape 5 "Rectangle 55" ---
$tr = $slide.Shapes.Item(5).TextFrame.TextRange
$tr.Characters(1, 2).Text = '{' + $q
$tr.Characters(12, 10).Text = $q + ': [0, 1],'
$tr.Characters(23, 2).Text = ' ' + $q
# Closing quote for "newTraceIndes" was missing before; the fix adds it.
$tr.Characters(38, 9).Text = $q + ': [1, 0]}'

# --- Shape 6 "Rectangle 56" ---
$tr = $slide.Shapes.Item(6).TextFrame.TextRange
$tr.Characters(1, 2).Text = '{' + $q
$tr.Characters(12, 7).Text = $q + ': [1]}'

# --- Shape 7 "Rectangle 57" ---
$tr = $slide.Shapes.Item(7).TextFrame.TextRange
$tr.Characters(1, 10).Text = '{' + $q + 'layout' + $q + ':'
$tr.Characters(12, 7).Text = '     {' + $q
$tr.Characters(30, 11).Text = $q + ': [-1, 3],'
$tr.Characters(42, 7).Text = '      ' + $q
$tr.Characters(60, 12).Text = $q + ': [-3, 5]}}'
